# Week 2 presentation: bump the "Codefeed" subtitle on slide 1 up to 36pt.
# (PowerPoint auto-adds <a:normAutofit/> to the shape's bodyPr the first
# time a run's font size is changed on a placeholder that doesn't already
# specify an autofit behavior.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Subtitle 2")

$tr = $shp.TextFrame.TextRange
$tr.Font.Size = 36
